# Update Name of Algo
# Applies the targeted numeric corrections to result_data_RandomForest.xlsx
# as described by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C3"  = -11.3523
    "C21" = -12.93980000000001
    "C23" = -12.3867
    "C25" = -13.4603
    "E27" = 16.71499999999999
    "E31" = 16.7864
    "E39" = 15.98099999999999
    "E48" = 17.4714
    "E51" = 16.8171
    "E52" = 16.6475
    "C53" = -10.5491
    "E55" = 16.54250000000001
    "E56" = 16.3644
    "C57" = -14.09019999999999
    "E57" = 16.81959999999999
    "C59" = -12.8065
    "C69" = -10.9834
    "E73" = 17.2436
    "C79" = -10.86510000000001
    "C83" = -13.7741
    "E89" = 17.35650000000002
    "E90" = 16.53209999999999
    "C93" = -11.3694
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
